$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C6").Value = 60000429
$ws.Range("C7").Value = 60000429
$ws.Range("D8").Value = 294
$ws.Range("D9").Value = 294
$ws.Range("C10").Value = 60000429
$ws.Range("D11").Value = 294
$ws.Range("D12").Value = 294
$ws.Range("D13").Value = 294
$ws.Range("D14").Value = 294
$ws.Range("C15").Value = 60000430
$ws.Range("D16").Value = 295
